# Supplementary Procedure.docx - apply the recorded edits.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the three split runs of the first OSF hyperlink
#    ("https://osf.io/n" + "k" + "r9g/") into a single run of text,
#    without disturbing the Hyperlink run style.
# ---------------------------------------------------------------------------
$h = $d.Hyperlinks.Item(1)
$hFind = $h.Range
$hFind.Find.Execute("https://osf.io/nkr9g/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hStart = $hFind.Start
$hText = $hFind.Text
$hLen = $hText.Length
# Round-trip through a differently-sized placeholder so the engine actually
# performs a text replace (identical text would be a no-op) and coalesces
# the three backing runs into one.
$placeholderHyperlink = "*" * ($hLen + 1)
$hRange = $d.Range($hStart, $hStart + $hLen)
$hRange.Text = $placeholderHyperlink
$hRange2 = $d.Range($hStart, $hStart + $placeholderHyperlink.Length)
$hRange2.Text = $hText

# ---------------------------------------------------------------------------
# 2) Merge the runs making up the "Study 7 ..." bullet into a single run,
#    then mark that whole paragraph with a tracked yellow-highlight format
#    change (author "Ian Hussey").
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    if ($text.Length -ge 7 -and $text.Substring(0, 7) -eq "Study 7") {
        $target = $para
        break
    }
}
$pRange = $target.Range
$bodyStart = $pRange.Start
$bodyEnd = $pRange.End - 1          # exclude the trailing paragraph mark
$bodyRange = $d.Range($bodyStart, $bodyEnd)
$bodyText = $bodyRange.Text
$placeholder = "*" * $bodyText.Length
$bodyRange.Text = $placeholder
$bodyRange2 = $d.Range($bodyStart, $bodyStart + $placeholder.Length)
$bodyRange2.Text = $bodyText

$word.UserName = "Ian Hussey"
$d.TrackRevisions = $true
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    if ($text.Length -ge 7 -and $text.Substring(0, 7) -eq "Study 7") {
        $target2 = $para
        break
    }
}
$target2.Range.HighlightColorIndex = 7
$d.TrackRevisions = $false

# ---------------------------------------------------------------------------
# 3) "D1 algorithm" -> "D2 algorithm" recorded as a tracked insertion of "2"
#    immediately followed by a tracked deletion of "1" (author "Ian Hussey").
# ---------------------------------------------------------------------------
$word.UserName = "Ian Hussey"
$d.TrackRevisions = $true

$findRange = $d.Content
$findRange.Find.Execute("the D1 algorithm", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchStart = $findRange.Start
$matchText = $findRange.Text
$dIndex = $matchText.IndexOf("D1")
$posAfterD = $matchStart + $dIndex + 1

$insPoint = $d.Range($posAfterD, $posAfterD)
$insPoint.InsertAfter("2")

$posOne = $posAfterD + 1
$oneRange = $d.Range($posOne, $posOne + 1)
$oneRange.Delete()

$d.TrackRevisions = $false

# ---------------------------------------------------------------------------
# 4) Move the (hidden) _GoBack bookmark from the end of the "raw data"
#    paragraph to the blank paragraph that follows the "Study 7..." bullet.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$blankTarget = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    if ($text.Length -ge 7 -and $text.Substring(0, 7) -eq "Study 7") {
        $blankTarget = $d.Paragraphs($i + 1)
        break
    }
}
$d.Bookmarks.Add("_GoBack", $blankTarget.Range)
